$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-43: refreshed Price (D) / Volume(1h) (E) figures ---
# Leading "'" forces text storage (these look numeric but are plain
# text cells in the source sheet, e.g. "26.854.04", "1.636.23").
$ws.Range("D2").Value = '''26.854.04'
$ws.Range("E2").Value = '''  +0.18%  '
$ws.Range("D3").Value = '''1.636.23'
$ws.Range("E3").Value = '''  -0.14%  '
$ws.Range("E4").Value = '''  -0.47%  '
$ws.Range("D5").Value = '''216.69'
$ws.Range("E5").Value = '''  -0.99%  '
$ws.Range("E6").Value = '''  +1.94%  '
$ws.Range("E7").Value = '''  -0.43%  '
$ws.Range("E8").Value = '''  +1.67%  '
$ws.Range("D9").Value = '''0.0624'
$ws.Range("E9").Value = '''  +0.24%  '
$ws.Range("D10").Value = '''19.95'
$ws.Range("E11").Value = '''  -0.03%  '
$ws.Range("D12").Value = '''1.866.43'
$ws.Range("E12").Value = '''  -0.10%  '
$ws.Range("D13").Value = '''1.641.23'
$ws.Range("E13").Value = '''  -0.01%  '
$ws.Range("E14").Value = '''  -0.80%  '
$ws.Range("E15").Value = '''  +0.38%  '
$ws.Range("E16").Value = '''  +2.58%  '
$ws.Range("D17").Value = '''26.857.61'
$ws.Range("E17").Value = '''  +0.20%  '
$ws.Range("E18").Value = '''  -0.55%  '
$ws.Range("D19").Value = '''219.65'
$ws.Range("E19").Value = '''  +1.63%  '
$ws.Range("E20").Value = '''  -0.47%  '
$ws.Range("D21").Value = '''6.76'
$ws.Range("E21").Value = '''  +2.89%  '
$ws.Range("E22").Value = '''  +0.66%  '
$ws.Range("D23").Value = '''2.45'
$ws.Range("E23").Value = '''  +3.87%  '
$ws.Range("E24").Value = '''  -0.21%  '
$ws.Range("D25").Value = '''147.04'
$ws.Range("E25").Value = '''  -0.09%  '
$ws.Range("E26").Value = '''  -0.35%  '
$ws.Range("E27").Value = '''  +3.32%  '
$ws.Range("E28").Value = '''  +0.45%  '
$ws.Range("D29").Value = '''15.79'
$ws.Range("E29").Value = '''  +0.31%  '
$ws.Range("E30").Value = '''  -0.40%  '
$ws.Range("E31").Value = '''  -1.43%  '
$ws.Range("E32").Value = '''  -1.34%  '
$ws.Range("E33").Value = '''  +0.67%  '
$ws.Range("E34").Value = '''  +0.61%  '
$ws.Range("D35").Value = '''1.254.95'
$ws.Range("E35").Value = '''  -0.63%  '
$ws.Range("E36").Value = '''  -0.22%  '
$ws.Range("E37").Value = '''  +1.75%  '
$ws.Range("E38").Value = '''  +0.31%  '
$ws.Range("D39").Value = '''0.832'
$ws.Range("E40").Value = '''  -0.41%  '
$ws.Range("D41").Value = '''0.811'
$ws.Range("E41").Value = '''  +0.54%  '
$ws.Range("E42").Value = '''  +1.22%  '
$ws.Range("D43").Value = '''1.776.16'
$ws.Range("E43").Value = '''  -0.18%  '

# --- Rows 44-51: coin list reshuffled (names/links/prices/volumes) ---
# Only cells whose value actually differs from the diff are rewritten,
# e.g. row 47 (RenderToken) keeps its original Price (D47) untouched.
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").Value = '''61.80'
$ws.Range("E44").Value = '''  +2.65%  '
$ws.Range("B45").Value = 'MXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D45").Value = '''2.10'
$ws.Range("E45").Value = '''  -1.59%  '
$ws.Range("D46").Value = '''91.56'
$ws.Range("E46").Value = '''  -0.71%  '
$ws.Range("E47").Value = '''  +0.01%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '''0.0513'
$ws.Range("E48").Value = '''  -0.55%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''7.67'
$ws.Range("E49").Value = '''  +1.93%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '''0.0960'
$ws.Range("E50").Value = '''  -0.53%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = '''0.404'
$ws.Range("E51").Value = '''  -0.21%  '
